$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.918.49'
$ws.Range("E2").Value = '  +1.39%  '
$ws.Range("D3").Value = '2.307.84'
$ws.Range("E3").Value = '  -0.13%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '542.20'
$ws.Range("E5").Value = '  +0.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.37'
$ws.Range("E6").Value = '  -2.34%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("E8").Value = '  -2.47%  '
$ws.Range("D9").Value = '2.306.11'
$ws.Range("E9").Value = '  -0.15%  '
$ws.Range("E10").Value = '  -0.28%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.55'
$ws.Range("E11").Value = '  +2.23%  '
$ws.Range("E12").Value = '  -0.24%  '
$ws.Range("E13").Value = '  -0.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.35'
$ws.Range("E14").Value = '  -2.17%  '
$ws.Range("D15").Value = '59.869.55'
$ws.Range("E15").Value = '  +1.46%  '
$ws.Range("D16").Value = '2.716.71'
$ws.Range("E16").Value = '  -0.41%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000131'
$ws.Range("E17").Value = '  -1.11%  '
$ws.Range("D18").Value = '2.309.26'
$ws.Range("E18").Value = '  +0.19%  '
$ws.Range("E19").Value = '  -1.27%  '
$ws.Range("E20").Value = '  -2.24%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '311.69'
$ws.Range("E21").Value = '  -0.37%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.55'
$ws.Range("E22").Value = '  -0.49%  '
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.66'
$ws.Range("E24").Value = '  +1.53%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.170'
$ws.Range("E25").Value = '  -1.76%  '
$ws.Range("E26").Value = '  -0.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.74'
$ws.Range("E27").Value = '  -2.60%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.34'
$ws.Range("E28").Value = '  +3.42%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '171.33'
$ws.Range("E29").Value = '  +0.77%  '
$ws.Range("E30").Value = '  -0.89%  '
$ws.Range("E31").Value = '  +0.07%  '
$ws.Range("D32").Value = '0.0₃0725'
$ws.Range("E32").Value = '  -1.82%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.82'
$ws.Range("E33").Value = '  -1.31%  '
$ws.Range("E34").Value = '  +3.08%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.378'
$ws.Range("E35").Value = '  -1.43%  '
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("E37").Value = '  -0.95%  '
$ws.Range("E38").Value = '  -0.05%  '
$ws.Range("E39").Value = '  -1.90%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '316.17'
$ws.Range("E40").Value = '  +2.21%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '38.07'
$ws.Range("E41").Value = '  -0.92%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.51'
$ws.Range("E42").Value = '  -0.71%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '135.88'
$ws.Range("E43").Value = '  -3.51%  '
$ws.Range("E44").Value = '  -1.00%  '
$ws.Range("E45").Value = '  -2.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '18.74'
$ws.Range("E47").Value = '  +1.87%  '
$ws.Range("E48").Value = '  -1.05%  '
$ws.Range("D49").Value = '0.0₆0224'
$ws.Range("E49").Value = '  +21.75%  '
$ws.Range("E50").Value = '  +0.54%  '
$ws.Range("E51").Value = '  +0.13%  '
